# add check balance, check mini statement, check custom statement
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the new account opening deposit amount into D2 (opening balance for account 90216/Current)
$ws.Range("D2").Value = 84915

# Remove the extra sample rows (rows 3 & 4 fully cleared incl. style), keeping only the D column placeholder
$ws.Range("A3").Clear()
$ws.Range("B3").Clear()
$ws.Range("C3").Clear()

$ws.Range("A4").Clear()
$ws.Range("B4").Clear()
$ws.Range("C4").Clear()

# Row 5: clear A5/B5 entirely, but keep C5's formatting (clear only its contents)
$ws.Range("A5").Clear()
$ws.Range("B5").Clear()
$ws.Range("C5").ClearContents()

# Update the current selection to highlight the remaining data rows
$ws.Range("A3:D5").Select()
